$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NIG(0.6942633242570948, 0.3981461132546792, 1.4961720406795351, 2.5994876987185243)"
$ws.Range("C2").Value = "NIG(0.6069255158594491, 0.4003641969867814, 5.958451715171844, 4.685976978937511)"
$ws.Range("D2").Value = "F(206.4510646425552, 16.277850638133582, -9.296039296200984, 12.509690660363582)"
$ws.Range("E2").Value = "JSU(-1.3749131636640661, 1.28512617685212, 2.118928461850451, 5.611909937533855)"
